$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H1").Value = 15.0
$ws.Range("I1").Value = 32.0
